$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ C = -0.3623658873974311;   E = 0.1825419310453658 }
    3  = @{ C = -0.009261555895478946; E = 0.1145211022186787 }
    4  = @{ C = -1.404263945418582;    E = -0.807808220045203 }
    5  = @{ C = 1.692932643509848;     E = 0.6262577107155831 }
    6  = @{ C = 1.020829760720643;     E = 1.148272834981245 }
    7  = @{ C = 0.6772121200332215;    E = 1.258913537332895 }
    8  = @{ C = 1.019715257608911;     E = 0.9536145745415947 }
    9  = @{ C = 2.173959184500385;     E = 1.566646323486043 }
    10 = @{ C = 1.707434489469994;     E = 1.30258347990615 }
    11 = @{ C = 1.456988786619839;     E = 1.842797144428188 }
    12 = @{ C = 1.241332692055597;     E = 1.58004210678635 }
    13 = @{ C = 1.592885137608979;     E = 1.604795846351514 }
    14 = @{ C = -2.015335584265165;    E = -1.215549235925828 }
    15 = @{ C = -3.579597300369253;    E = -1.403103901755631 }
    16 = @{ C = 4.461954539041502;     E = 0.7797949948739058 }
    17 = @{ C = -1.305206755692701;    E = 0.5821000732047832 }
    18 = @{ C = 0.0845726262934221;    E = 0.1341520870597357 }
    19 = @{ C = 0.9724700385226326;    E = 0.6236501628417823 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
